# Add a new "PF/1.0.2" row to the meta-sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.2"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# The new row keeps the default/"Normal" formatting (no inherited column style)
$ws.Range("A3:D3").Style = "Normal"
